# Update "想去人数" (want-to-go count) figures in column F across the
# relevant worksheets, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 129
$wsExhibit.Range("F3").Value = 2148
$wsExhibit.Range("F4").Value = 30
$wsExhibit.Range("F5").Value = 11241
$wsExhibit.Range("F10").Value = 11165
$wsExhibit.Range("F11").Value = 451
$wsExhibit.Range("F15").Value = 5589
$wsExhibit.Range("F17").Value = 3449

# --- Sheet "演出" (performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 567

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 129
$wsAll.Range("F3").Value = 2148
$wsAll.Range("F4").Value = 567
$wsAll.Range("F5").Value = 30
$wsAll.Range("F7").Value = 11241
$wsAll.Range("F12").Value = 11165
$wsAll.Range("F13").Value = 451
$wsAll.Range("F17").Value = 5589
$wsAll.Range("F19").Value = 3449
